# no-op for now
$p = $ppt.ActivePresentation
Write-Output ("Slides.Count=" + $p.Slides.Count)
